$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45945
$ws.Range("B2").Value = 116.04
$ws.Range("C2").Value = 108.26
$ws.Range("D2").Value = 105.45
$ws.Range("E2").Value = 103.77
$ws.Range("F2").Value = 103.77
$ws.Range("G2").Value = 106.96
$ws.Range("H2").Value = 113.9
$ws.Range("I2").Value = 142.4
$ws.Range("J2").Value = 158.17
$ws.Range("K2").Value = 121.45
$ws.Range("L2").Value = 100.75
$ws.Range("M2").Value = 85.28
$ws.Range("N2").Value = 82
$ws.Range("O2").Value = 81.26000000000001
$ws.Range("P2").Value = 77.43000000000001
$ws.Range("Q2").Value = 81.81999999999999
$ws.Range("R2").Value = 82.23
$ws.Range("S2").Value = 100.34
$ws.Range("T2").Value = 116.88
$ws.Range("U2").Value = 156.52
$ws.Range("V2").Value = 163.07
$ws.Range("W2").Value = 148.73
$ws.Range("X2").Value = 120.25
$ws.Range("Y2").Value = 112.67
$ws.Range("Z2").Value = 112.06
$ws.Range("AB2").Value = 136.18
$ws.Range("AD2").Value = 155.9
$ws.Range("AF2").Value = 139.81
